# Macroferia Regional de Talca - Repollo
# Insert two new daily price records at the top of the data block (rows 268-269),
# pushing the existing rows down by two (old row N becomes new row N+2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 268..364 down by two rows, making room for the two new records.
$ws.Rows("268:269").Insert()

# New row 268: Crespo record / Primera
$ws.Range("A268").Value = 5
$ws.Range("B268").Value = "Macroferia Regional de Talca"
$ws.Range("C268").Value = "Maule"
$ws.Range("D268").Value = 44809
$ws.Range("E268").Value = 7
$ws.Range("F268").Value = 100112006
$ws.Range("G268").Value = "Repollo"
$ws.Range("H268").Value = "Crespo record"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 2000
$ws.Range("K268").Value = 1300
$ws.Range("L268").Value = 1300
$ws.Range("M268").Value = 1300
$ws.Range("N268").Value = "`$/unidad"
$ws.Range("O268").Value = "Región del Maule"
$ws.Range("P268").Value = 1300
$ws.Range("Q268").Value = 1
$ws.Range("R268").Value = "Hortaliza"

# New row 269: Crespo record / Segunda
$ws.Range("A269").Value = 5
$ws.Range("B269").Value = "Macroferia Regional de Talca"
$ws.Range("C269").Value = "Maule"
$ws.Range("D269").Value = 44809
$ws.Range("E269").Value = 7
$ws.Range("F269").Value = 100112006
$ws.Range("G269").Value = "Repollo"
$ws.Range("H269").Value = "Crespo record"
$ws.Range("I269").Value = "Segunda"
$ws.Range("J269").Value = 2000
$ws.Range("K269").Value = 1000
$ws.Range("L269").Value = 1000
$ws.Range("M269").Value = 1000
$ws.Range("N269").Value = "`$/unidad"
$ws.Range("O269").Value = "Región del Maule"
$ws.Range("P269").Value = 1000
$ws.Range("Q269").Value = 1
$ws.Range("R269").Value = "Hortaliza"

# Match the date cell number format used by the rest of column D.
$ws.Range("D268:D269").NumberFormat = $ws.Range("D270").NumberFormat
